$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (want-to-go count)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 981
$wsExpo.Range("F7").Value = 7611
$wsExpo.Range("F13").Value = 4730
$wsExpo.Range("F17").Value = 4998
$wsExpo.Range("F23").Value = 242
$wsExpo.Range("F26").Value = 8619
$wsExpo.Range("F31").Value = 65
$wsExpo.Range("F33").Value = 66
$wsExpo.Range("F37").Value = 1772
$wsExpo.Range("F39").Value = 1063

# Sheet "全部类型" (All types) - same events mirrored at different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 981
$wsAll.Range("F9").Value = 7611
$wsAll.Range("F15").Value = 4730
$wsAll.Range("F19").Value = 4998
$wsAll.Range("F25").Value = 242
$wsAll.Range("F29").Value = 8619
$wsAll.Range("F34").Value = 65
$wsAll.Range("F35").Value = 66
$wsAll.Range("F38").Value = 1772
$wsAll.Range("F40").Value = 1063
